$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 14 for lab task "9.2", shifting existing
# rows 14-35 down to 15-36 (formulas/refs shift automatically).
$ws.Rows.Item(14).Insert()

# Label the newly inserted row and its "chapters completed" formula.
$ws.Range("A14").Value = "9.2"
$ws.Range("B14").Formula = "=SUM(C14:N14)"

# Link lab task "9.2" (row 14) to chapter 6 (column F).
$ws.Range("F14").Value = 1

# Link lab task "6.1" (row 5) to chapter 14 (column M).
$ws.Range("M5").Value = 1

# Restore the active selection to M11.
[void]$ws.Range("M11").Select()
